$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1758
$ws.Range("F3").Value = 10276
$ws.Range("F6").Value = 602
$ws.Range("F8").Value = 1660
$ws.Range("F9").Value = 196
$ws.Range("F10").Value = 413
$ws.Range("F12").Value = 226
$ws.Range("F13").Value = 499
$ws.Range("F14").Value = 1186
$ws.Range("F15").Value = 135
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 107
$ws.Range("F20").Value = 373
$ws.Range("F21").Value = 18
$ws.Range("F22").Value = 340
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 108
$ws.Range("F25").Value = 1176
$ws.Range("F26").Value = 27
$ws.Range("F27").Value = 42
$ws.Range("F29").Value = 249
$ws.Range("F31").Value = 556
$ws.Range("F32").Value = 237
$ws.Range("F33").Value = 380
$ws.Range("F35").Value = 708
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 777
$ws.Range("F40").Value = 824
$ws.Range("F41").Value = 712
$ws.Range("F42").Value = 393
$ws.Range("F43").Value = 349
$ws.Range("F44").Value = 31
$ws.Range("F45").Value = 363
$ws.Range("F46").Value = 83
$ws.Range("F47").Value = 357

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 7
$ws.Range("F13").Value = 1
$ws.Range("F18").Value = 1117
$ws.Range("F20").Value = 1687
$ws.Range("F21").Value = 1688
$ws.Range("F22").Value = 1109
$ws.Range("F23").Value = 338
$ws.Range("F25").Value = 82
$ws.Range("F26").Value = 13
$ws.Range("F27").Value = 26
$ws.Range("F29").Value = 372
$ws.Range("F32").Value = 218
$ws.Range("F40").Value = 133
$ws.Range("F44").Value = 76

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 818
$ws.Range("F5").Value = 203
$ws.Range("F6").Value = 2546
$ws.Range("F7").Value = 4171
$ws.Range("F8").Value = 71
$ws.Range("F10").Value = 364
$ws.Range("F11").Value = 178
$ws.Range("F12").Value = 217

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1758
$ws.Range("F3").Value = 818
$ws.Range("F5").Value = 10276
$ws.Range("F6").Value = 203
$ws.Range("F7").Value = 4171
$ws.Range("F8").Value = 71
$ws.Range("F9").Value = 364
$ws.Range("F10").Value = 365
$ws.Range("F11").Value = 602
$ws.Range("F12").Value = 1660
$ws.Range("F13").Value = 196
$ws.Range("F15").Value = 7
$ws.Range("F17").Value = 135
$ws.Range("F23").Value = 107
$ws.Range("F24").Value = 1117
$ws.Range("F25").Value = 373
$ws.Range("F26").Value = 340
$ws.Range("F27").Value = 1109
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 26
$ws.Range("F30").Value = 42
$ws.Range("F31").Value = 249
$ws.Range("F32").Value = 372
$ws.Range("F33").Value = 556
$ws.Range("F35").Value = 380
$ws.Range("F38").Value = 218
$ws.Range("F39").Value = 777
$ws.Range("F41").Value = 824
$ws.Range("F42").Value = 712
$ws.Range("F43").Value = 393
$ws.Range("F44").Value = 349
$ws.Range("F45").Value = 31
$ws.Range("F46").Value = 133
$ws.Range("F47").Value = 363
$ws.Range("F48").Value = 357
$ws.Range("F49").Value = 76
